$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DPI number
$ws.Range("A2").Value = 3306466721202

# Update birth date text (shared string is edited in place)
$ws.Range("B2").Value = "18/09/1997"

# Match the pasted-in font formatting applied to the birth date cell
$ws.Range("B2").Font.Name = "Roboto"
$ws.Range("B2").Font.Size = 8
$ws.Range("B2").Font.Color = 2367776

# Update the active selection to the birth date cell
$ws.Range("B2").Select() | Out-Null
